$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6564
$ws.Range("F5").Value = 109
$ws.Range("F6").Value = 609
$ws.Range("F7").Value = 212
$ws.Range("F8").Value = 40
$ws.Range("F9").Value = 790
$ws.Range("F10").Value = 1255
$ws.Range("F14").Value = 496
$ws.Range("F16").Value = 1036
$ws.Range("F17").Value = 1447
$ws.Range("F18").Value = 701
$ws.Range("F19").Value = 424
$ws.Range("F20").Value = 424
$ws.Range("F21").Value = 90
$ws.Range("F22").Value = 1090
$ws.Range("F23").Value = 200
$ws.Range("F24").Value = 2277
$ws.Range("F26").Value = 144
$ws.Range("F29").Value = 3706
$ws.Range("F31").Value = 677

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 31
$ws.Range("F6").Value = 728
$ws.Range("F11").Value = 130
$ws.Range("F19").Value = 4104
$ws.Range("F24").Value = 211
$ws.Range("F26").Value = 99
$ws.Range("F29").Value = 37

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1597
$ws.Range("F7").Value = 135
$ws.Range("F8").Value = 893

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1597
$ws.Range("F7").Value = 135
$ws.Range("F8").Value = 6564
$ws.Range("F9").Value = 31
$ws.Range("F11").Value = 728
$ws.Range("F12").Value = 212
$ws.Range("F13").Value = 40
$ws.Range("F14").Value = 790
$ws.Range("F17").Value = 130
$ws.Range("F18").Value = 130
$ws.Range("F22").Value = 1255
$ws.Range("F25").Value = 496
$ws.Range("F29").Value = 1036
$ws.Range("F30").Value = 1447
$ws.Range("F32").Value = 701
$ws.Range("F33").Value = 424
$ws.Range("F34").Value = 424
$ws.Range("F35").Value = 90
$ws.Range("F37").Value = 211
$ws.Range("F39").Value = 99
$ws.Range("F44").Value = 144
$ws.Range("F47").Value = 3706
$ws.Range("F51").Value = 677
